$wb = $excel.ActiveWorkbook

# --- Sheet "Customer": update group/category values and selection ---
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Range("B2").Value = "Gold"
$wsCustomer.Range("A2").Value = "Thanh_06092023"
$wsCustomer.Range("D6").Select()

# --- Sheet "Edit Customer": update group/category values and selection ---
$wsEditCustomer = $wb.Worksheets.Item("Edit Customer")
$wsEditCustomer.Range("A2").Value = "Thanh_07092023"
$wsEditCustomer.Range("B2").Value = "VIP"
$wsEditCustomer.Range("D9").Select()

# --- Sheet "Project": update the customer group reference and selection ---
$wsProject = $wb.Worksheets.Item("Project")
$wsProject.Range("B2").Value = "Thanh_06092023"
$wsProject.Range("C12").Select()

# --- Sheet "Edit Project": selection only ---
$wsEditProject = $wb.Worksheets.Item("Edit Project")
$wsEditProject.Range("C10").Select()
